# Userstories sheet cleanup: correct Sprint/Status ("Done") info for a
# few rows that were left stale or blank, and drop the now-unused "..."
# placeholder string. Also refresh the saved scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 12: Sprint was wrong (29) and Status still showed the "..." placeholder.
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = "Done"

# Row 18: Sprint/Status had never been filled in.
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = "Done"

# Row 19: Sprint/Status had never been filled in.
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = "Done"

# Refresh the view that gets persisted with the sheet: scroll back near the
# top of the data and leave the active selection on the Status cell we just
# finished checking.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.ScrollRow = 1
} catch {
    # Scroll position is cosmetic only; ignore if unsupported.
}
$ws.Range("G17").Select()
